$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("homePage")
$ws2 = $wb.Worksheets.Item("upgradeNow")

# --- homePage sheet changes ---
# B2 numeric value 10 -> 50
$ws1.Range("B2").Value = 50

# Clear out the contents of A6/B6 (previously "timeout" / 3), keep formatting
$ws1.Range("A6:B6").ClearContents()

# --- upgradeNow sheet changes ---
# A3 action changes from "visit" to "navigate"
$ws2.Range("A3").Value = "navigate"

# --- Active tab / selection changes ---
# upgradeNow becomes the selected/active tab, with A3 selected
# homePage keeps its previous selection of A4, but is no longer the active tab
$ws1.Range("A4").Select()
$ws2.Activate()
$ws2.Range("A3").Select()
